# Updates the "合肥-漫展信息" workbook (gh-pages data refresh at 456a3b4):
#  - refreshes "想去人数" (column F) counts on a handful of existing rows
#  - inserts a new event row ("合肥·第二届TH动漫游戏嘉年华", 2024-07-28) in
#    date order on both the "展览" sheet and the "全部类型" sheet
#  - bumps the "第七届环形宇宙动漫游戏嘉年华" and "银魂主题派对only2.0" counts
#    that shift down to make room for the new row

$wb = $excel.ActiveWorkbook

# (sheet index, row the new event should be inserted at -- i.e. the row
#  currently holding "第七届环形宇宙动漫游戏嘉年华")
$sheetPlans = @(
    @{ Index = 1; InsertBefore = 21 },
    @{ Index = 4; InsertBefore = 21 }
)

# Plain "想去人数" (column F) refreshes, keyed by row number - identical on
# both sheets since the tables share the same leading rows.
$fUpdates = @{
    2  = 1572
    3  = 8904
    4  = 101
    6  = 669
    7  = 331
    9  = 38
    10 = 54
    11 = 3772
    13 = 374
    14 = 98
    15 = 4085
    16 = 4
    20 = 233
}

# Writes a literal string into a cell without Excel's "looks like a date"
# auto-conversion kicking in (which would turn "2024-07-28" into a date
# serial + date number format). Entering it as a `="..."` text formula and
# immediately flattening it to a static value keeps the cell a plain string
# with the default (unformatted) style, matching the source data.
function Set-LiteralText($cell, $text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

foreach ($plan in $sheetPlans) {
    $ws = $wb.Worksheets.Item($plan.Index)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value2 = $fUpdates[$row]
    }

    $insertRow = $plan.InsertBefore

    # Shift everything from $insertRow down by one row, then fill the freed
    # row with the new event, preserving column-A's bold/bordered style.
    $ws.Rows.Item($insertRow).Insert(-4121, 0)

    $ws.Cells.Item($insertRow - 1, 1).Copy()
    $ws.Cells.Item($insertRow, 1).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    $ws.Cells.Item($insertRow, 1).Value2 = $ws.Cells.Item($insertRow - 1, 1).Value2 + 1
    Set-LiteralText $ws.Cells.Item($insertRow, 2) "2024-07-28"
    $ws.Cells.Item($insertRow, 3).Value2 = "合肥·第二届TH动漫游戏嘉年华"
    $ws.Cells.Item($insertRow, 4).Value2 = "田埠西路199号 吉祥如意宴会楼蜀山店"
    $ws.Cells.Item($insertRow, 5).Value2 = "2024.07.28 09:30-07.28 17:00"
    $ws.Cells.Item($insertRow, 6).Value2 = 2
    $ws.Cells.Item($insertRow, 7).Value2 = 55
    $ws.Cells.Item($insertRow, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=87447"
    $ws.Cells.Item($insertRow, 9).Value2 = "//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"

    # Renumber every row below the inserted one (column A is a simple
    # 0-based running index = row - 1).
    $dims = $ws.UsedRange
    $lastRow = $dims.Row + $dims.Rows.Count - 1
    for ($r = $insertRow + 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $r - 1
    }

    # The rows pushed down by the insert need their refreshed "想去人数"
    # values (everything else about them is unchanged).
    for ($r = $insertRow + 1; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($name -eq "合肥·第七届环形宇宙动漫游戏嘉年华") {
            $ws.Cells.Item($r, 6).Value2 = 2569
        } elseif ($name -eq "合肥·银魂主题派对only2.0") {
            $ws.Cells.Item($r, 6).Value2 = 97
        }
    }
}
